# Generate Report for Handoff
# Updates the localization-status workbook to reflect that the
# f1607909-...-md file's zh-cn/de-de handoff has just been (re)generated:
#   - Status moves from "Handed back: in sync with en-US" to "Ready for handoff"
#     for both language rows (row 2 = c68f627d file, row 3 = f1607909 file)
#   - The f1607909 row's "Latest Handoff Datetime" is refreshed
#   - The f1607909 row's "Error Detail" now explains the handback file is stale
#   - The "Error Detail" column is widened to fit the new message

$wb = $excel.ActiveWorkbook

$statusReady = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fc875d6c33694db627d53acbf97734b336bfc936/e2e/f1607909-33e9-459d-84e0-46fd80efbca0.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dabd8479449c97aa7d8be5412f0851a33410675a/e2e/f1607909-33e9-459d-84e0-46fd80efbca0.md."

# --- Overview sheet: row 3 corresponds to f1607909-...-md ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusReady
$wsOverview.Range("F3").Value = $statusReady
$wsOverview.Range("G3").Value = "2016-08-19 04:43:18"

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $statusReady
$wsZh.Range("C3").Value = $statusReady
$wsZh.Range("H3").Value = "2016-08-19 04:43:14"
$wsZh.Range("P3").Value = $errorDetail
$wsZh.Range("P1").ColumnWidth = 39.166666666666664

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $statusReady
$wsDe.Range("C3").Value = $statusReady
$wsDe.Range("H3").Value = "2016-08-19 04:43:18"
$wsDe.Range("P3").Value = $errorDetail
$wsDe.Range("P1").ColumnWidth = 39.166666666666664

Write-Output "Updated handoff status for f1607909-33e9-459d-84e0-46fd80efbca0"
